$wb = $excel.ActiveWorkbook

$portugal = $wb.Worksheets.Item("Portugal")
$portugal.Activate()
$portugal.Cells.Select()

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Copy($null, $last)

$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"
$slovakia.Activate()

$slovakia.Range("B4").Value = "NGC-2930/T3242/T3241/T3243"
$slovakia.Range("B4").Style = "Normal"

$slovakia.Rows("3:5").RowHeight = 14.4
$slovakia.Rows("7:7").RowHeight = 43.2

$slovakia.Range("B4").Select()

Write-Host "done"
